$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.4532
$ws.Range("A3").Value = -21.46230000000002
$ws.Range("C5").Value = -12.911
$ws.Range("A14").Value = -20.44129999999998
$ws.Range("A21").Value = -21.23030000000001
$ws.Range("A23").Value = -21.87010000000003
$ws.Range("A25").Value = -22.38870000000002
